$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row: value was blank, now set to the text "true".
# A leading apostrophe forces Excel to store it as literal text instead
# of auto-converting to the Boolean TRUE.
$ws.Range("B7").Value = "'true"

# "Date" row: value updated to new timestamp
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
